# Phase 3 RAD Non-UI Test Cases and Data
# Append the new RAD test-data rows (TaxType coverage for additional
# Existing Liability w/Notice Number and New Tax Return Amount Due cases)
# to the bottom of the EmailNoMatch RAD data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New rows 34-41: PaymentType = "Existing Liability w/Notice Number" ---
$existingLiabilityTaxTypes = @(
    "Admissions and Amusement Tax",
    "Estate Tax",
    "Motor Fuel Tax",
    "Slots License Fee",
    "Tobacco Tax",
    "Transportation Network Services",
    "Unclaimed Property",
    "IFTA Tax"
)

$r = 34
foreach ($taxType in $existingLiabilityTaxTypes) {
    $ws.Range("C$r").Value = "Y"
    $ws.Range("D$r").Value = "Existing Liability w/Notice Number"
    $ws.Range("E$r").Value = $taxType
    $r = $r + 1
}

# --- New rows 42-54: PaymentType = "New Tax Return Amount Due" ---
$newTaxReturnTaxTypes = @(
    "Admissions and Amusement Tax",
    "Alcohol Tax",
    "Bay Restoration Fee",
    "Corporate Income Tax",
    "Estate Tax",
    "Motor Fuel Tax",
    "Sales and Use Tax",
    "Slots License Fee",
    "Tire Recycling Fee",
    "Tobacco Tax",
    "Transportation Network Services",
    "Unclaimed Property",
    "Withholding Tax"
)

foreach ($taxType in $newTaxReturnTaxTypes) {
    $ws.Range("C$r").Value = "Y"
    $ws.Range("D$r").Value = "New Tax Return Amount Due"
    $ws.Range("E$r").Value = $taxType
    $r = $r + 1
}

# --- Update the view so the newly-added rows are visible/selected ---
$ws.Application.ActiveWindow.ScrollRow = 22
[void]$ws.Range("C28:C54").Select()

Write-Host "Added RAD test rows 34-54"
